$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-05 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-06 Wednesday", 2)

$d.Content.Find.Execute("620÷5=124, 0", $true, $false, $false, $false, $false, $true, 1, $false, "428÷9=47, 5", 2)
$d.Content.Find.Execute("341÷3=113, 2", $true, $false, $false, $false, $false, $true, 1, $false, "561÷6=93, 3", 2)
$d.Content.Find.Execute("791÷5=158, 1", $true, $false, $false, $false, $false, $true, 1, $false, "788÷5=157, 3", 2)
$d.Content.Find.Execute("497÷9=55, 2", $true, $false, $false, $false, $false, $true, 1, $false, "687÷9=76, 3", 2)
$d.Content.Find.Execute("156÷3=52, 0", $true, $false, $false, $false, $false, $true, 1, $false, "879÷6=146, 3", 2)

$d.Content.Find.Execute("731÷9=81, 2", $true, $false, $false, $false, $false, $true, 1, $false, "355÷5=71, 0", 2)
$d.Content.Find.Execute("414÷4=103, 2", $true, $false, $false, $false, $false, $true, 1, $false, "593÷3=197, 2", 2)
$d.Content.Find.Execute("618÷4=154, 2", $true, $false, $false, $false, $false, $true, 1, $false, "601÷7=85, 6", 2)
$d.Content.Find.Execute("130÷3=43, 1", $true, $false, $false, $false, $false, $true, 1, $false, "926÷7=132, 2", 2)
$d.Content.Find.Execute("409÷9=45, 4", $true, $false, $false, $false, $false, $true, 1, $false, "621÷6=103, 3", 2)

$d.Content.Find.Execute("131÷7=18, 5", $true, $false, $false, $false, $false, $true, 1, $false, "225÷2=112, 1", 2)
$d.Content.Find.Execute("441÷6=73, 3", $true, $false, $false, $false, $false, $true, 1, $false, "250÷7=35, 5", 2)
$d.Content.Find.Execute("436÷8=54, 4", $true, $false, $false, $false, $false, $true, 1, $false, "387÷7=55, 2", 2)
$d.Content.Find.Execute("471÷4=117, 3", $true, $false, $false, $false, $false, $true, 1, $false, "107÷2=53, 1", 2)
$d.Content.Find.Execute("745÷4=186, 1", $true, $false, $false, $false, $false, $true, 1, $false, "336÷3=112, 0", 2)

$d.Content.Find.Execute("919÷9=102, 1", $true, $false, $false, $false, $false, $true, 1, $false, "517÷5=103, 2", 2)
$d.Content.Find.Execute("221÷7=31, 4", $true, $false, $false, $false, $false, $true, 1, $false, "337÷4=84, 1", 2)
$d.Content.Find.Execute("696÷9=77, 3", $true, $false, $false, $false, $false, $true, 1, $false, "632÷4=158, 0", 2)
$d.Content.Find.Execute("687÷8=85, 7", $true, $false, $false, $false, $false, $true, 1, $false, "125÷3=41, 2", 2)
$d.Content.Find.Execute("513÷9=57, 0", $true, $false, $false, $false, $false, $true, 1, $false, "487÷5=97, 2", 2)

$d.Content.Find.Execute("721÷6=120, 1", $true, $false, $false, $false, $false, $true, 1, $false, "769÷5=153, 4", 2)
$d.Content.Find.Execute("490÷9=54, 4", $true, $false, $false, $false, $false, $true, 1, $false, "993÷9=110, 3", 2)
$d.Content.Find.Execute("881÷2=440, 1", $true, $false, $false, $false, $false, $true, 1, $false, "918÷6=153, 0", 2)
$d.Content.Find.Execute("203÷8=25, 3", $true, $false, $false, $false, $false, $true, 1, $false, "689÷8=86, 1", 2)
$d.Content.Find.Execute("839÷3=279, 2", $true, $false, $false, $false, $false, $true, 1, $false, "576÷6=96, 0", 2)
